$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A264").Value = "IMX-USD"
$ws.Range("A265").Value = "MNT-USD"
$ws.Range("A266").Value = "TAO-USD"
